# 11-6-1-1.xlsx — update responsible-organization / contact-person / website
# metadata (Section 2) and move the active selection to B4.
#
# Shared-string table note: the workbook stores column-B text as shared
# strings. Overwriting a cell's .Value drops the old string (if no longer
# referenced) and appends the new one at the end of the table, which is
# exactly how the authored diff re-numbers the <si> entries — so the other
# B-cells (whose *text* is unchanged) automatically pick up their shifted
# indices with no further action needed here. The order below (B7, then
# B10, then B6) reproduces the same append order as the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Contact person (B7): Керималиева Н.К. -> Мамбеталиев Т.А.
$ws.Range("B7").Value = "Мамбеталиев Т.А."

# Organization website (B10): www.stat.kg -> www.stat.gov.kg
$ws.Range("B10").Value = "www.stat.gov.kg "

# Organization (B6): statistics-of-sustainable-development dept -> digital
# development & statistics-of-sustainable-development dept
$ws.Range("B6").Value = "Национальный статистический комитет КР (Управление цифрового развития и статистики устойчивого развития)"

# Move/save the active cell selection from B2 to B4
$ws.Range("B4").Select()
